# Apply the edits described by the diff:
#  - Number the existing "Inspiration" questions 1)-4)
#  - Append four new numbered questions (5-8) plus a trailing blank
#    paragraph, right after question 4) and before the two blank
#    paragraphs that precede the "OBJECTIVES" heading.

$d = $word.ActiveDocument

function Prefix-Paragraph($searchText, $prefix) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $insertPoint = $d.Range($rng.Start, $rng.Start)
        $insertPoint.InsertBefore($prefix)
    }
}

# 1) - 4): prepend numbering to the existing questions.
Prefix-Paragraph "Is the number of petitions with Data Engineer" "1)"
Prefix-Paragraph "Which part of the US has the most Hardware Engineer jobs?" "2)"
Prefix-Paragraph "Which industry has the most number of Data Scientist positions?" "3)"
Prefix-Paragraph "Which employers" "4)"

# Locate question 4's paragraph so the new questions can be appended
# immediately after it (same run formatting: sz=32 / szCs=32).
$rng4 = $d.Content
$rng4.Find.Execute("4)Which employers(CORPORATES AND COMPANIES) file the most petitions each year?", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$para4 = $rng4.Paragraphs(1)
$insertRange = $para4.Range
$insertRange.Collapse(0)

$cr = [char]13
$block = "5)Which year(2011-2016) had most H1B petitions approved and for which Job Positions ?" + $cr `
       + "6)Predictive modelling on for which Job titles H1B Visa petitions of foreign Nationals will be approved?" + $cr `
       + "7)Predicting  the no of H1B applications going to increases or decrease for a particular Job position in future?" + $cr `
       + "8)Predicting the Prevailing wages of particular Job Titles for future?" + $cr

$insertRange.InsertParagraphAfter()
$insertRange.InsertAfter($block)
